$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at W (pushes the existing "nota_iniciativa" column,
# and everything to its right, one column to the right -> becomes X).
$ws.Columns("W:W").Insert()

# New column W is the "l1" repeat-group placeholder column, mirroring "c1" (column V).
$ws.Range("W1").Value = "l1"

# Fill the new column's data rows (2-84) with 0, matching the sibling c1/nota_iniciativa columns.
$ws.Range("W2:W84").Value = 0
